$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data rows (rows 2-5) below the "Saldo inicial" row,
# then lay the data back out the way add_data_to_excel does: a blank
# spacer row, the column headers, and a single new income/expense entry.
$ws.Range("A2:C5").Clear()

# Row 2: blank spacer row
$ws.Range("A2").Value = " "

# Row 3: header row
$ws.Range("A3").Value = "Ingresos"
$ws.Range("B3").Value = "Gastos"
$ws.Range("C3").Value = "Fecha"

# Row 4: newly added data entry. Set the date format before the value so
# Excel doesn't stamp its own default short-date format on the cell first.
$ws.Range("A4").Value = " "
$ws.Range("B4").Value = 60
$ws.Range("C4").NumberFormat = "yyyy-mm-dd"
$ws.Range("C4").Value = Get-Date -Year 2020 -Month 9 -Day 9 -Hour 0 -Minute 0 -Second 0
